$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new record right above the current row 121, pushing the
# existing rows 121-197 down to 122-198 (dimension grows from R197 to R198).
$ws.Rows("121:121").Insert()

# Populate the newly inserted row 121 with the new "Albahaca" record.
$ws.Cells.Item(121, 1).Value = 4
$ws.Cells.Item(121, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(121, 3).Value = "Los Lagos"
$ws.Cells.Item(121, 4).Value = 45090
$ws.Cells.Item(121, 5).Value = 10
$ws.Cells.Item(121, 6).Value = 100112052
$ws.Cells.Item(121, 7).Value = "Albahaca"
$ws.Cells.Item(121, 8).Value = "Sin especificar"
$ws.Cells.Item(121, 9).Value = "Primera"
$ws.Cells.Item(121, 10).Value = 80
$ws.Cells.Item(121, 11).Value = 5000
$ws.Cells.Item(121, 12).Value = 5000
$ws.Cells.Item(121, 13).Value = 5000
$ws.Cells.Item(121, 14).Value = "$/paquete"
$ws.Cells.Item(121, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(121, 16).Value = 5000
$ws.Cells.Item(121, 17).Value = 1
$ws.Cells.Item(121, 18).Value = "Hortaliza"
